# Dep Ed Closures.xlsx - update content from "Monday 14 September" snapshot
# to "Tuesday 15 September" snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "On this page" summary cell (row 7): day/date/time refresh.
$ws.Cells.Item(7, 1).Value = "On this pageCurrent school and early childhood service, TAFE closures and relocations:Bus service cancellations or alterationsCurrent school" + [char]160 + "and early childhood service," + [char]160 + "TAFE closures and relocations for Tuesday 15" + [char]160 + "September," + [char]160 + "(as at 9:50am, 15" + [char]160 + "September)South-Eastern Victoria RegionEarly childhood services"

# 2) Insert a new early-childhood-service closure entry before the current
#    row 104 ("li: St Andrews Christian College Outside School Hours Care WANTIRNA SOUTH").
$ws.Rows.Item(104).Insert()
$ws.Cells.Item(104, 1).Value = "li: Ruyton Early Learning Kindergarten & Pre Prep KEW"

# That insertion shifts "li: Hopetoun Child Care Service HOPETOUN" down from row 180
# to row 181; it is no longer present in the updated list, so remove it.
$ws.Rows.Item(181).Delete()

# 3) "li: Wydinia Kindergarten & Early Learning Centre COLAC" (originally row 202,
#    now row 202 again after the +1/-1 shifts above) is no longer present either.
$ws.Rows.Item(202).Delete()

# 4) Add the newly closed school to the South-Western Victoria region school list.
$ws.Cells.Item(205, 1).Value = "Al Taqwa College, TRUGANINAOur lady of the Southern Cross, MANOR LAKESParkville College (Malmsbury campus), PARKVILLETAFE"
